$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "NA" value from C119 down to the new row 120, and append a new
# data row (120) that mirrors row 119's other values with the next date.
# A bare "'" forces an empty TEXT value (rather than clearing the cell to
# blank), matching the other empty-page-number cells in the column; reset
# the style afterwards so the quote-prefix formatting doesn't stick.
$ws.Cells.Item(119, 3).Value = "'"
$ws.Cells.Item(119, 3).Style = "Normal"

# Use a leading apostrophe so the date-shaped string is stored as literal
# text (matching the rest of the sheet) instead of being parsed into a
# date serial number, then reset the style so no extra formatting sticks.
$ws.Cells.Item(120, 1).Value = "'2025-05-28"
$ws.Cells.Item(120, 1).Style = "Normal"
$ws.Cells.Item(120, 2).Value = "Rien ne nous concerne aujourd'hui !"
$ws.Cells.Item(120, 3).Value = "NA"
$ws.Cells.Item(120, 4).Value = 1
